$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dynamics")

$path1 = "k:\github\digitalmodel\tests\modules\orcaflex\orcaflex_post_process\orcaflex_test1.sim"
$path2 = "k:\github\digitalmodel\tests\modules\orcaflex\orcaflex_post_process\orcaflex_test2.sim"

$ws.Range("B2").Value = $path1
$ws.Range("D2").Value = $path1

$ws.Range("B3").Value = $path2
$ws.Range("D3").Value = $path2
